# Updates cryptos price list (columns B-E) to the latest scraped values.
# D-column price cells are text (e.g. '307.70', '42.500.94') so each is written
# with a leading apostrophe to keep Excel from re-interpreting it as a number
# (which would otherwise drop trailing zeros / thousands separators), then the
# cell style is reset to Normal so no stray 'Text' number-format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'42.500.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.91%  "
# Row 3
$ws.Range("D3").Value = "'2.292.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.39%  "
# Row 4
$ws.Range("E4").Value = "  -0.15%  "
# Row 5
$ws.Range("D5").Value = "'157.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15,607.31%  "
# Row 6
$ws.Range("D6").Value = "'307.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "
# Row 7
$ws.Range("D7").Value = "'95.91"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.39%  "
# Row 8
$ws.Range("E8").Value = "  +0.42%  "
# Row 9
$ws.Range("E9").Value = "  -0.14%  "
# Row 10
$ws.Range("D10").Value = "'0.497"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.64%  "
# Row 11
$ws.Range("D11").Value = "'36.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.74%  "
# Row 12
$ws.Range("E12").Value = "  +1.19%  "
# Row 13
$ws.Range("E13").Value = "  -1.96%  "
# Row 14
$ws.Range("E14").Value = "  +2.77%  "
# Row 15
$ws.Range("D15").Value = "'2.646.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "
# Row 16
$ws.Range("E16").Value = "  +2.87%  "
# Row 17
$ws.Range("D17").Value = "'2.300.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.12%  "
# Row 18
$ws.Range("D18").Value = "'0.803"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.01%  "
# Row 19
$ws.Range("D19").Value = "'42.401.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "
# Row 20
$ws.Range("D20").Value = "'12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.83%  "
# Row 21
$ws.Range("D21").Value = "'0.0₃0919"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.11%  "
# Row 22
$ws.Range("E22").Value = "  +2.35%  "
# Row 23
$ws.Range("D23").Value = "'68.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.32%  "
# Row 24
$ws.Range("D24").Value = "'243.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.33%  "
# Row 25
$ws.Range("E25").Value = "  +1.45%  "
# Row 26
$ws.Range("E26").Value = "  +2.66%  "
# Row 27
$ws.Range("E27").Value = "  -0.08%  "
# Row 28
$ws.Range("D28").Value = "'24.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.44%  "
# Row 29
$ws.Range("D29").Value = "'35.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.35%  "
# Row 30
$ws.Range("D30").Value = "'9.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.24%  "
# Row 31
$ws.Range("E31").Value = "  +1.33%  "
# Row 32
$ws.Range("D32").Value = "'161.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "
# Row 33
$ws.Range("E33").Value = "  +4.21%  "
# Row 35
$ws.Range("D35").Value = "'0.0757"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.01%  "
# Row 36
$ws.Range("E36").Value = "  +3.28%  "
# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.109"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.93%  "
# Row 38
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'17.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.57%  "
# Row 39
$ws.Range("E39").Value = "  -0.26%  "
# Row 40
$ws.Range("D40").Value = "'1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.98%  "
# Row 41
$ws.Range("E41").Value = "  +0.10%  "
# Row 42
$ws.Range("E42").Value = "  +6.58%  "
# Row 43
$ws.Range("D43").Value = "'2.011.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.33%  "
# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.29%  "
# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'2.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.99%  "
# Row 46
$ws.Range("E46").Value = "  +3.05%  "
# Row 47
$ws.Range("D47").Value = "'3.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.07%  "
# Row 48
$ws.Range("D48").Value = "'10.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "
# Row 49
$ws.Range("D49").Value = "'53.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.82%  "
# Row 50
$ws.Range("E50").Value = "  +2.84%  "
# Row 51
$ws.Range("D51").Value = "'73.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
